# Update "want to go" counts (column F) on several rows across sheets,
# as published in the gh-pages regenerated data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 2042
$ws1.Range("F23").Value = 2845
$ws1.Range("F26").Value = 3209
$ws1.Range("F27").Value = 663
$ws1.Range("F29").Value = 235
$ws1.Range("F33").Value = 701
$ws1.Range("F34").Value = 675

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 79
$ws2.Range("F16").Value = 148
$ws2.Range("F21").Value = 189

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 391

# Sheet "全部类型" (All types) - aggregated view of the above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F15").Value = 391
$ws4.Range("F16").Value = 2042
$ws4.Range("F17").Value = 2042
$ws4.Range("F26").Value = 79
$ws4.Range("F38").Value = 2845
$ws4.Range("F40").Value = 3209
$ws4.Range("F41").Value = 663
$ws4.Range("F43").Value = 235
$ws4.Range("F51").Value = 701
$ws4.Range("F52").Value = 675
